# Weekly price-list update: a new observation is inserted at row 158
# (Feria Lagunitas de Puerto Montt - Betarraga), pushing the previously
# existing rows 158-188 down to 159-189.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 158; Excel shifts rows 158:188 down
# to 159:189 and the sheet's used range grows to A1:R189 automatically.
$ws.Rows("158:158").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A158").Value = 4
$ws.Range("B158").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C158").Value = "Los Lagos"
$ws.Range("D158").Value = 44474
$ws.Range("E158").Value = 10
$ws.Range("F158").Value = 100114014
$ws.Range("G158").Value = "Betarraga"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 1400
$ws.Range("K158").Value = 1000
$ws.Range("L158").Value = 1000
$ws.Range("M158").Value = 1000
$ws.Range("N158").Value = "$/paquete 5 unidades"
$ws.Range("O158").Value = "Región del Maule"
$ws.Range("P158").Value = 200
$ws.Range("Q158").Value = 5
$ws.Range("R158").Value = "Hortaliza"
